# Motor Controller BOM update
# - Swap the 2200 uF 63V electrolytic cap for the 50V variant (shorter part,
#   allows a shorter controller design).
# - Swap the MAX17640CATA (adjustable 0.9V) buck regulator for the
#   MAX17640BATA (fixed 5V) variant.
# - Leave the cursor on C5 like the author did before saving.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: 2200 uF Aluminum Electrolytic Capacitor ---------------------
$ws.Range("B4").Value = "1189-4388-ND"
$ws.Range("C4").Value = "2200 µF 50 V Aluminum Electrolytic Capacitors Radial, Can - 10000 Hrs @ 105°C"
$ws.Range("E4").Value = 2.84

$ws.Range("H4").Value = "https://www.digikey.com/en/products/detail/rubycon/50ZLJ2200M18X25/10437363"
$hlB4 = $ws.Range("H4").Hyperlinks.Item(1)
$hlB4.Address = "https://www.digikey.com/en/products/detail/rubycon/50ZLJ2200M18X25/10437363"

# --- Row 12: MAX17640 Buck Switching Regulator ---------------------------
$ws.Range("B12").Value = "175-MAX17640BATA+-ND"
$ws.Range("C12").Value = "Buck Switching Regulator IC Positive Fixed 5V 1 Output 400mA 8-WFDFN"

$ws.Range("H12").Value = "https://www.digikey.com/en/products/detail/maxim-integrated/MAX17640BATA/14287859?s=N4IgTCBcDaIIwHYCsBaAsgQQBqIGwBYAGAIQwBUMBqFAOQBEQBdAXyA"
$hlB12 = $ws.Range("H12").Hyperlinks.Item(1)
$hlB12.Address = "https://www.digikey.com/en/products/detail/maxim-integrated/MAX17640BATA/14287859?s=N4IgTCBcDaIIwHYCsBaAsgQQBqIGwBYAGAIQwBUMBqFAOQBEQBdAXyA"

# --- Cosmetic: restore the author's last selection before saving --------
$ws.Range("C5").Select()

$wb.Save()
